$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (closest achievable via ColumnWidth COM property,
# which snaps to whole-pixel boundaries under the hood)
$ws.Columns.Item(1).ColumnWidth = 15.666666666666668
$ws.Columns.Item(2).ColumnWidth = 14.833333333333332

# Update cell values
$ws.Range("A1").Value = -0.089808712187860351
$ws.Range("B1").Value = 0.089424213299338362
$ws.Range("A2").Value = -0.043313021271684349
$ws.Range("B2").Value = 0.04173774064917346
$ws.Range("A3").Value = 0.11574122265354347
$ws.Range("B3").Value = -0.11627530436470224
$ws.Range("A4").Value = -0.16772105924575342
$ws.Range("B4").Value = 0.166850556560032
$ws.Range("A5").Value = -0.1608505567775822
$ws.Range("B5").Value = 0.15910791519564071
$ws.Range("A6").Value = -0.04632993517199635
$ws.Range("B6").Value = 0.046303677617765437
$ws.Range("A7").Value = -0.026303677884495613
$ws.Range("B7").Value = 0.026269841172741337
$ws.Range("A8").Value = -0.0062698414404200875
$ws.Range("B8").Value = 0.0062480313558639367
$ws.Range("A9").Value = -0.00024803158280395365
$ws.Range("B9").Value = 0.00022368070155387443
$ws.Range("A10").Value = 0.030889997995231511
$ws.Range("B10").Value = -0.030929744897328248
$ws.Range("A11").Value = 0.035429744674956964
$ws.Range("B11").Value = -0.035523348377608954
$ws.Range("A12").Value = -0.020153479237583038
$ws.Range("B12").Value = 0.020108747892487067
$ws.Range("A13").Value = -0.014108748120237991
$ws.Range("B13").Value = 0.01410120770610046
$ws.Range("A14").Value = -0.0021012079518873961
$ws.Range("B14").Value = 0.0021008462141880457
$ws.Range("A15").Value = -0.021051136153229066
$ws.Range("B15").Value = 0.021026624367998181
$ws.Range("A16").Value = -0.015026624597005878
$ws.Range("B16").Value = 0.015004398516011541
$ws.Range("A17").Value = -0.0090043987460628472
$ws.Range("B17").Value = 0.0089999997608440907
$ws.Range("A18").Value = -0.036111192254317359
$ws.Range("B18").Value = 0.036096856237779917
$ws.Range("A19").Value = -0.027096856461699126
$ws.Range("B19").Value = 0.027013922214804253
$ws.Range("A20").Value = -0.018013922440813346
$ws.Range("B20").Value = 0.018004272070724525
$ws.Range("A21").Value = -0.0090042722970613553
$ws.Range("B21").Value = 0.0089999997733896109
$ws.Range("A22").Value = -0.14847896352769752
$ws.Range("B22").Value = 0.14745265200102509
$ws.Range("A23").Value = -0.084631998528297281
$ws.Range("B23").Value = 0.084126397641222539
$ws.Range("A24").Value = -0.042126397970635665
$ws.Range("B24").Value = 0.041999999668608012
$ws.Range("A25").Value = -0.094884222148870379
$ws.Range("B25").Value = 0.094642097493387922
$ws.Range("A26").Value = -0.088642097719418445
$ws.Range("B26").Value = 0.088330956266929661
$ws.Range("A27").Value = -0.082330956494184093
$ws.Range("B27").Value = 0.081269973104475035
$ws.Range("A28").Value = -0.075269973336794749
$ws.Range("B28").Value = 0.074534729298054714
$ws.Range("A29").Value = -0.062534729551693147
$ws.Range("B29").Value = 0.062169736491538785
$ws.Range("A30").Value = -0.042169736770436117
$ws.Range("B30").Value = 0.042020282892756722
$ws.Range("A31").Value = -0.027020283159471248
$ws.Range("B31").Value = 0.027000910629546127
$ws.Range("A32").Value = -0.0060009109144054662
$ws.Range("B32").Value = 0.0059999997591599907
